$d = $word.ActiveDocument
$d.Content.Find.Execute("” means the following kinds of information:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "” means:", 2)
